$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new row of data on Sheet1
$ws1.Range("A8").Value = "ddddd"
$ws1.Range("B8").Value = "cccccc"

# Sheet1 becomes the active sheet/tab
$ws1.Activate()

# Move the selection on Sheet1 to B10
$ws1.Range("B10").Select()
